# Auto-verifiable edit script for AUS.xlsx VAR/ARIMA bugfix
$wb = $excel.ActiveWorkbook

# --- Sheet 1: y_fitted_on_begin_2016 ---
# Prepend a new oldest-year row (1991) and recompute every subsequent value.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Insert()
$ws1.Rows.Item(2).ClearFormats()
$ws1.Range("A2").Value = 1991
$ws1.Range("B2").Value = 13.18009964968262
$ws1.Range("A3").Value = 1992
$ws1.Range("B3").Value = 12.72425138211406
$ws1.Range("A4").Value = 1993
$ws1.Range("B4").Value = 12.65785286803129
$ws1.Range("A5").Value = 1994
$ws1.Range("B5").Value = 12.77440581109808
$ws1.Range("A6").Value = 1995
$ws1.Range("B6").Value = 12.85449441158461
$ws1.Range("A7").Value = 1996
$ws1.Range("B7").Value = 12.66909203392301
$ws1.Range("A8").Value = 1997
$ws1.Range("B8").Value = 12.43188907280924
$ws1.Range("A9").Value = 1998
$ws1.Range("B9").Value = 12.13498291317227
$ws1.Range("A10").Value = 1999
$ws1.Range("B10").Value = 11.92065881213938
$ws1.Range("A11").Value = 2000
$ws1.Range("B11").Value = 11.72709196822052
$ws1.Range("A12").Value = 2001
$ws1.Range("B12").Value = 11.3094661429657
$ws1.Range("A13").Value = 2002
$ws1.Range("B13").Value = 11.01718630668362
$ws1.Range("A14").Value = 2003
$ws1.Range("B14").Value = 10.60391915963347
$ws1.Range("A15").Value = 2004
$ws1.Range("B15").Value = 10.54219003490809
$ws1.Range("A16").Value = 2005
$ws1.Range("B16").Value = 10.29180908672768
$ws1.Range("A17").Value = 2006
$ws1.Range("B17").Value = 10.05106056763595
$ws1.Range("A18").Value = 2007
$ws1.Range("B18").Value = 9.497316915870762
$ws1.Range("A19").Value = 2008
$ws1.Range("B19").Value = 9.137611773638881
$ws1.Range("A20").Value = 2009
$ws1.Range("B20").Value = 8.405494326485879
$ws1.Range("A21").Value = 2010
$ws1.Range("B21").Value = 7.707667213300669
$ws1.Range("A22").Value = 2011
$ws1.Range("B22").Value = 7.298249155968027
$ws1.Range("A23").Value = 2012
$ws1.Range("B23").Value = 7.280665715942083
$ws1.Range("A24").Value = 2013
$ws1.Range("B24").Value = 6.626158854802686
$ws1.Range("A25").Value = 2014
$ws1.Range("B25").Value = 6.415249349161683
$ws1.Range("A26").Value = 2015
$ws1.Range("B26").Value = 6.236882760797225
$ws1.Range("A27").Value = 2016
$ws1.Range("B27").Value = 6.140361965445657

# --- Sheet 2: y_pred_on_2017_2021 ---
# Years unchanged, only y_value recomputed.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = 6.113437268770937
$ws2.Range("B3").Value = 6.129324017637515
$ws2.Range("B4").Value = 6.207925654185885
$ws2.Range("B5").Value = 6.31145923050218
$ws2.Range("B6").Value = 6.490043515426692

# --- Sheet 3: y_fitted_on_begin_2021 ---
# Drop the final (2021 old) row and recompute every value; years shift to 1991-2021.
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(33).Delete()
$ws3.Range("A2").Value = 1991
$ws3.Range("B2").Value = 13.23884058761216
$ws3.Range("A3").Value = 1992
$ws3.Range("B3").Value = 12.68472044400386
$ws3.Range("A4").Value = 1993
$ws3.Range("B4").Value = 12.63593205229932
$ws3.Range("A5").Value = 1994
$ws3.Range("B5").Value = 12.76920385377928
$ws3.Range("A6").Value = 1995
$ws3.Range("B6").Value = 12.88042078829129
$ws3.Range("A7").Value = 1996
$ws3.Range("B7").Value = 12.70185643657713
$ws3.Range("A8").Value = 1997
$ws3.Range("B8").Value = 12.48872221247919
$ws3.Range("A9").Value = 1998
$ws3.Range("B9").Value = 12.10046809698646
$ws3.Range("A10").Value = 1999
$ws3.Range("B10").Value = 11.97945967695825
$ws3.Range("A11").Value = 2000
$ws3.Range("B11").Value = 11.67210176314776
$ws3.Range("A12").Value = 2001
$ws3.Range("B12").Value = 11.32340652335218
$ws3.Range("A13").Value = 2002
$ws3.Range("B13").Value = 10.96709215899097
$ws3.Range("A14").Value = 2003
$ws3.Range("B14").Value = 10.56153424571948
$ws3.Range("A15").Value = 2004
$ws3.Range("B15").Value = 10.53490958799802
$ws3.Range("A16").Value = 2005
$ws3.Range("B16").Value = 10.32478181272105
$ws3.Range("A17").Value = 2006
$ws3.Range("B17").Value = 10.03341459593763
$ws3.Range("A18").Value = 2007
$ws3.Range("B18").Value = 9.418313522834472
$ws3.Range("A19").Value = 2008
$ws3.Range("B19").Value = 9.107426972157848
$ws3.Range("A20").Value = 2009
$ws3.Range("B20").Value = 8.353862654043201
$ws3.Range("A21").Value = 2010
$ws3.Range("B21").Value = 7.721057356777965
$ws3.Range("A22").Value = 2011
$ws3.Range("B22").Value = 7.289872646153971
$ws3.Range("A23").Value = 2012
$ws3.Range("B23").Value = 7.189347402481694
$ws3.Range("A24").Value = 2013
$ws3.Range("B24").Value = 6.874165774715665
$ws3.Range("A25").Value = 2014
$ws3.Range("B25").Value = 6.557044399436457
$ws3.Range("A26").Value = 2015
$ws3.Range("B26").Value = 6.160628200216358
$ws3.Range("A27").Value = 2016
$ws3.Range("B27").Value = 6.027454955369224
$ws3.Range("A28").Value = 2017
$ws3.Range("B28").Value = 5.971869049168824
$ws3.Range("A29").Value = 2018
$ws3.Range("B29").Value = 5.737572591407001
$ws3.Range("A30").Value = 2019
$ws3.Range("B30").Value = 5.581456198314782
$ws3.Range("A31").Value = 2020
$ws3.Range("B31").Value = 5.526337618525091
$ws3.Range("A32").Value = 2021
$ws3.Range("B32").Value = 5.630699411273427

# --- Sheet 4: y_pred_on_2022_2026 ---
# Years unchanged, only y_value recomputed.
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = 6.234246727068481
$ws4.Range("B3").Value = 6.790643801957671
$ws4.Range("B4").Value = 7.086344500236706
$ws4.Range("B5").Value = 7.34735383675627
$ws4.Range("B6").Value = 7.656039958812215

Write-Output "edit applied"
